$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 250
$ws.Range("C2").Value = 320
$ws.Range("D2").Value = 22
$ws.Range("E2").Value = 3065.056341566094
$ws.Range("F2").Value = 6.796725879902986

# Row 3
$ws.Range("B3").Value = 221.5
$ws.Range("C3").Value = 320
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 3065.056341566094
$ws.Range("F3").Value = 6.796725879902986

# Row 4
$ws.Range("B4").Value = 221.5
$ws.Range("C4").Value = 174.8941768640577
$ws.Range("E4").Value = 2817.006492904334
$ws.Range("F4").Value = 7.289754939244287

# Row 5
$ws.Range("B5").Value = 28.50000000000002
$ws.Range("C5").Value = 319.9999999999999
$ws.Range("D5").Value = 22
$ws.Range("E5").Value = 3065.056341566094
$ws.Range("F5").Value = 6.796725879902986

# Row 6
$ws.Range("B6").Value = 28.50000000000002
$ws.Range("E6").Value = 2313.856579806076
$ws.Range("F6").Value = 7.393562710918144
$ws.Range("H6").Value = 89.07955179903253

# Row 7
$ws.Range("B7").Value = 28.50000000000002
$ws.Range("E7").Value = 173.8398003706387

# Row 8
$ws.Range("B8").Value = 28.50000000000002
$ws.Range("C8").Value = 41.53552743031628
$ws.Range("E8").Value = 174.1652036954729
$ws.Range("F8").Value = 0.5927497189592826

# Row 9
$ws.Range("B9").Value = 19.00950000000001
$ws.Range("C9").Value = 41.53552743031628
$ws.Range("E9").Value = 174.1652036954729
$ws.Range("F9").Value = 0.5927497189592826

# Row 10
$ws.Range("B10").Value = 9.490500000000004
$ws.Range("C10").Value = 41.53552743031628
$ws.Range("E10").Value = 174.1652036954729
$ws.Range("F10").Value = 0.5927497189592826

# Row 11
$ws.Range("B11").Value = 221.5
$ws.Range("C11").Value = 174.8941768640577
$ws.Range("E11").Value = 2817.006492904334
$ws.Range("F11").Value = 7.289754939244287

# Row 12
$ws.Range("B12").Value = 215.5195
$ws.Range("C12").Value = 174.8941768640577
$ws.Range("E12").Value = 2817.006492904334
$ws.Range("F12").Value = 7.289754939244287

# Row 13
$ws.Range("B13").Value = 225.01
$ws.Range("E13").Value = 2705.536401551786
$ws.Range("F13").Value = 7.025104808443722
$ws.Range("H13").Value = 99.49782569949188

# Row 14
$ws.Range("B14").Value = 225.01

# Row 15
$ws.Range("B15").Value = 5.9805
$ws.Range("C15").Value = 174.8941768640577
$ws.Range("E15").Value = 2817.006492904334
$ws.Range("F15").Value = 7.289754939244287

# Row 16
$ws.Range("B16").Value = 250
$ws.Range("C16").Value = 100.1968001138945
$ws.Range("E16").Value = 420.1073802389807
$ws.Range("F16").Value = 1.309317439218598

# Row 17
$ws.Range("B17").Value = 250
$ws.Range("C17").Value = 100.4932876541408
$ws.Range("D17").Value = 22
$ws.Range("E17").Value = 422.8205825802663
$ws.Range("F17").Value = 1.311135522056935
